$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Truncate/round the Ost (Q2) and Nord (R2) coordinate values to integers
$ws.Range("Q2").Value = 598458
$ws.Range("R2").Value = 6978200

# Clear the Starttid (Z2) and Sluttid (AB2) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
